$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44350
$ws.Range("K2").Value = 1800
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = 1900
$ws.Range("P2").Value = 633

# Row 3
$ws.Range("D3").Value = 44217
$ws.Range("J3").Value = 250
$ws.Range("K3").Value = 2900
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = 2950
$ws.Range("P3").Value = 983

# Row 4
$ws.Range("D4").Value = 44322
$ws.Range("J4").Value = 250
$ws.Range("K4").Value = 1400
$ws.Range("L4").Value = 1500
$ws.Range("M4").Value = 1450
$ws.Range("P4").Value = 483

# Row 5
$ws.Range("D5").Value = 44300
$ws.Range("J5").Value = 160
$ws.Range("K5").Value = 1000
$ws.Range("L5").Value = 1200
$ws.Range("M5").Value = 1100
$ws.Range("O5").Value = "Región de Coquimbo"
$ws.Range("P5").Value = 367

# Row 6
$ws.Range("D6").Value = 44545
$ws.Range("K6").Value = 2800
$ws.Range("L6").Value = 3000
$ws.Range("M6").Value = 2900
$ws.Range("P6").Value = 967

# Row 7
$ws.Range("D7").Value = 44600
$ws.Range("J7").Value = 320
$ws.Range("K7").Value = 1400
$ws.Range("L7").Value = 1500
$ws.Range("M7").Value = 1450
$ws.Range("P7").Value = 483

# Row 8
$ws.Range("D8").Value = 44320

# Row 9
$ws.Range("D9").Value = 44320

# Row 10
$ws.Range("D10").Value = 44579
$ws.Range("J10").Value = 300
$ws.Range("K10").Value = 3000
$ws.Range("L10").Value = 3500
$ws.Range("M10").Value = 3250
$ws.Range("P10").Value = 1083

# Row 11
$ws.Range("D11").Value = 44594
$ws.Range("I11").Value = "Segunda"
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 2000
$ws.Range("L11").Value = 2200
$ws.Range("M11").Value = 2100
$ws.Range("P11").Value = 700

# Row 12
$ws.Range("D12").Value = 44497
$ws.Range("K12").Value = 750
$ws.Range("L12").Value = 800
$ws.Range("M12").Value = 775
$ws.Range("P12").Value = 258

# Row 13
$ws.Range("D13").Value = 44308
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 270
$ws.Range("K13").Value = 1400
$ws.Range("L13").Value = 1500
$ws.Range("M13").Value = 1450
$ws.Range("P13").Value = 483

# Row 14
$ws.Range("D14").Value = 44495
$ws.Range("J14").Value = 300
$ws.Range("K14").Value = 950
$ws.Range("L14").Value = 1000
$ws.Range("M14").Value = 975
$ws.Range("P14").Value = 325

# Row 15
$ws.Range("D15").Value = 44448
$ws.Range("J15").Value = 200
$ws.Range("K15").Value = 1400
$ws.Range("L15").Value = 1500
$ws.Range("M15").Value = 1450
$ws.Range("P15").Value = 483

# Row 16
$ws.Range("D16").Value = 44448
$ws.Range("I16").Value = "Segunda"
$ws.Range("J16").Value = 200
$ws.Range("K16").Value = 1000
$ws.Range("L16").Value = 1200
$ws.Range("M16").Value = 1100
$ws.Range("P16").Value = 367

# Row 17
$ws.Range("D17").Value = 44327
$ws.Range("J17").Value = 200
$ws.Range("K17").Value = 1400
$ws.Range("L17").Value = 1500
$ws.Range("M17").Value = 1450
$ws.Range("P17").Value = 483

# Row 18
$ws.Range("D18").Value = 44327
$ws.Range("I18").Value = "Segunda"
$ws.Range("K18").Value = 1000
$ws.Range("L18").Value = 1200
$ws.Range("M18").Value = 1100
$ws.Range("P18").Value = 367

# Row 19
$ws.Range("D19").Value = 44257
$ws.Range("J19").Value = 1500
$ws.Range("K19").Value = 2800
$ws.Range("L19").Value = 3000
$ws.Range("M19").Value = 2900
$ws.Range("P19").Value = 967

# Row 20
$ws.Range("D20").Value = 44643
$ws.Range("J20").Value = 300
$ws.Range("K20").Value = 900
$ws.Range("L20").Value = 1000
$ws.Range("M20").Value = 950
$ws.Range("P20").Value = 317

# Row 21
$ws.Range("D21").Value = 44481
$ws.Range("K21").Value = 900
$ws.Range("L21").Value = 1000
$ws.Range("M21").Value = 950
$ws.Range("P21").Value = 317

# Row 22
$ws.Range("D22").Value = 44292
$ws.Range("J22").Value = 270
$ws.Range("M22").Value = 2450
$ws.Range("P22").Value = 817

# Row 23
$ws.Range("D23").Value = 44649
$ws.Range("J23").Value = 300
$ws.Range("K23").Value = 1800
$ws.Range("M23").Value = 1900
$ws.Range("P23").Value = 633

# Row 24
$ws.Range("D24").Value = 44435
$ws.Range("J24").Value = 270
$ws.Range("K24").Value = 1800
$ws.Range("L24").Value = 2000
$ws.Range("M24").Value = 1900
$ws.Range("P24").Value = 633

# Row 25
$ws.Range("D25").Value = 44483
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 300
$ws.Range("K25").Value = 1000
$ws.Range("L25").Value = 1200
$ws.Range("M25").Value = 1100
$ws.Range("P25").Value = 367

# Row 26
$ws.Range("D26").Value = 44460
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 300
$ws.Range("K26").Value = 950
$ws.Range("L26").Value = 1000
$ws.Range("M26").Value = 975
$ws.Range("P26").Value = 325

# Row 27
$ws.Range("D27").Value = 44341
$ws.Range("J27").Value = 300
$ws.Range("K27").Value = 1400
$ws.Range("L27").Value = 1500
$ws.Range("M27").Value = 1450
$ws.Range("P27").Value = 483

# Row 28
$ws.Range("D28").Value = 44525
$ws.Range("J28").Value = 250
$ws.Range("K28").Value = 1800
$ws.Range("L28").Value = 2000
$ws.Range("M28").Value = 1900
$ws.Range("P28").Value = 633

# Row 29
$ws.Range("D29").Value = 44586
$ws.Range("J29").Value = 250
$ws.Range("K29").Value = 2500
$ws.Range("M29").Value = 2750
$ws.Range("P29").Value = 917

# Row 30
$ws.Range("D30").Value = 44175
$ws.Range("K30").Value = 1800
$ws.Range("L30").Value = 2000
$ws.Range("M30").Value = 1900
$ws.Range("P30").Value = 633

# Row 31
$ws.Range("D31").Value = 44418
$ws.Range("K31").Value = 2400
$ws.Range("L31").Value = 2500
$ws.Range("M31").Value = 2450
$ws.Range("P31").Value = 817

# Row 32
$ws.Range("D32").Value = 44299
$ws.Range("J32").Value = 300
$ws.Range("K32").Value = 1400
$ws.Range("L32").Value = 1500
$ws.Range("M32").Value = 1450
$ws.Range("P32").Value = 483

# Row 33
$ws.Range("D33").Value = 44299
$ws.Range("I33").Value = "Segunda"
$ws.Range("J33").Value = 250

# Row 34
$ws.Range("D34").Value = 44642
$ws.Range("J34").Value = 250
$ws.Range("K34").Value = 1500
$ws.Range("M34").Value = 1750
$ws.Range("P34").Value = 583

# Row 35
$ws.Range("D35").Value = 44406
$ws.Range("J35").Value = 300
$ws.Range("K35").Value = 2800
$ws.Range("L35").Value = 3000
$ws.Range("M35").Value = 2900
$ws.Range("P35").Value = 967

# Row 36
$ws.Range("D36").Value = 44273
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 250
$ws.Range("K36").Value = 3800
$ws.Range("L36").Value = 4000
$ws.Range("M36").Value = 3900
$ws.Range("P36").Value = 1300

# Row 37
$ws.Range("D37").Value = 44356
$ws.Range("I37").Value = "Primera"
$ws.Range("K37").Value = 2400
$ws.Range("L37").Value = 2500
$ws.Range("M37").Value = 2450
$ws.Range("P37").Value = 817

# Row 38
$ws.Range("D38").Value = 44356
$ws.Range("I38").Value = "Segunda"
$ws.Range("J38").Value = 200
$ws.Range("K38").Value = 1800
$ws.Range("L38").Value = 2000
$ws.Range("M38").Value = 1900
$ws.Range("P38").Value = 633

# Row 39
$ws.Range("D39").Value = 44530
$ws.Range("J39").Value = 300
$ws.Range("K39").Value = 1900
$ws.Range("L39").Value = 2000
$ws.Range("M39").Value = 1950
$ws.Range("P39").Value = 650

# Row 40
$ws.Range("D40").Value = 44614
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 250
$ws.Range("K40").Value = 1500
$ws.Range("L40").Value = 2000
$ws.Range("M40").Value = 1750
$ws.Range("P40").Value = 583

# Row 41
$ws.Range("D41").Value = 44474
$ws.Range("J41").Value = 270
$ws.Range("K41").Value = 1000
$ws.Range("L41").Value = 1200
$ws.Range("M41").Value = 1100
$ws.Range("P41").Value = 367

# Row 42
$ws.Range("D42").Value = 44343
$ws.Range("I42").Value = "Primera"
$ws.Range("J42").Value = 150
$ws.Range("K42").Value = 1500
$ws.Range("L42").Value = 1500
$ws.Range("M42").Value = 1500
$ws.Range("P42").Value = 500

# Row 43
$ws.Range("D43").Value = 44343
$ws.Range("I43").Value = "Segunda"
$ws.Range("J43").Value = 150
$ws.Range("K43").Value = 1400
$ws.Range("L43").Value = 1400
$ws.Range("M43").Value = 1400
$ws.Range("P43").Value = 467

# Row 44
$ws.Range("D44").Value = 44277
$ws.Range("J44").Value = 250
$ws.Range("K44").Value = 1800
$ws.Range("L44").Value = 2000
$ws.Range("M44").Value = 1900
$ws.Range("P44").Value = 633

# Row 45
$ws.Range("D45").Value = 44335
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 250
$ws.Range("K45").Value = 1400
$ws.Range("L45").Value = 1500
$ws.Range("M45").Value = 1450
$ws.Range("P45").Value = 483

# Row 46
$ws.Range("D46").Value = 44392
$ws.Range("J46").Value = 200
$ws.Range("K46").Value = 3800
$ws.Range("L46").Value = 4000
$ws.Range("M46").Value = 3900
$ws.Range("P46").Value = 1300

# Row 47
$ws.Range("D47").Value = 44392
$ws.Range("I47").Value = "Segunda"
$ws.Range("J47").Value = 200
$ws.Range("K47").Value = 3200
$ws.Range("L47").Value = 3500
$ws.Range("M47").Value = 3350
$ws.Range("P47").Value = 1117

# Row 48
$ws.Range("D48").Value = 44565
$ws.Range("J48").Value = 250

# Row 49
$ws.Range("D49").Value = 44313
$ws.Range("J49").Value = 300
$ws.Range("K49").Value = 1300
$ws.Range("L49").Value = 1500
$ws.Range("M49").Value = 1400
$ws.Range("O49").Value = "Región de Arica y Parinacota"
$ws.Range("P49").Value = 467

# Row 50
$ws.Range("D50").Value = 44313
$ws.Range("I50").Value = "Segunda"
$ws.Range("J50").Value = 250
$ws.Range("K50").Value = 900
$ws.Range("L50").Value = 1000
$ws.Range("M50").Value = 950
$ws.Range("P50").Value = 317

# Row 51
$ws.Range("D51").Value = 44487
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 300
$ws.Range("K51").Value = 950
$ws.Range("L51").Value = 1000
$ws.Range("M51").Value = 975
$ws.Range("P51").Value = 325

# Row 52
$ws.Range("D52").Value = 44467
$ws.Range("K52").Value = 800
$ws.Range("L52").Value = 900
$ws.Range("M52").Value = 850
$ws.Range("P52").Value = 283

# Row 53
$ws.Range("D53").Value = 44558
$ws.Range("K53").Value = 3500
$ws.Range("L53").Value = 4000
$ws.Range("M53").Value = 3750
$ws.Range("P53").Value = 1250

# Row 54
$ws.Range("D54").Value = 44168
$ws.Range("J54").Value = 300
$ws.Range("K54").Value = 1800
$ws.Range("L54").Value = 2000
$ws.Range("M54").Value = 1900
$ws.Range("P54").Value = 633

# Row 55
$ws.Range("D55").Value = 44537
$ws.Range("J55").Value = 250
$ws.Range("K55").Value = 1400
$ws.Range("L55").Value = 1500
$ws.Range("M55").Value = 1450
$ws.Range("P55").Value = 483

# Row 56
$ws.Range("D56").Value = 44496
$ws.Range("J56").Value = 250
$ws.Range("K56").Value = 750
$ws.Range("L56").Value = 800
$ws.Range("M56").Value = 775
$ws.Range("P56").Value = 258

# Row 57
$ws.Range("D57").Value = 44412
$ws.Range("J57").Value = 300

# Row 58
$ws.Range("D58").Value = 44383
$ws.Range("I58").Value = "Segunda"
$ws.Range("J58").Value = 350
$ws.Range("K58").Value = 2800
$ws.Range("L58").Value = 3000
$ws.Range("M58").Value = 2886
$ws.Range("P58").Value = 962

# Row 59
$ws.Range("D59").Value = 44364
$ws.Range("J59").Value = 270
$ws.Range("K59").Value = 3400
$ws.Range("L59").Value = 3500
$ws.Range("M59").Value = 3450
$ws.Range("P59").Value = 1150

# Row 60
$ws.Range("D60").Value = 44635
$ws.Range("I60").Value = "Primera"
$ws.Range("J60").Value = 300
$ws.Range("K60").Value = 1800
$ws.Range("L60").Value = 2000
$ws.Range("M60").Value = 1900
$ws.Range("P60").Value = 633

# Row 61
$ws.Range("D61").Value = 44376
$ws.Range("J61").Value = 280
$ws.Range("K61").Value = 2400
$ws.Range("L61").Value = 2500
$ws.Range("M61").Value = 2436
$ws.Range("P61").Value = 812

# Row 62
$ws.Range("D62").Value = 44432
$ws.Range("J62").Value = 270
$ws.Range("K62").Value = 1800
$ws.Range("L62").Value = 2000
$ws.Range("M62").Value = 1900
$ws.Range("P62").Value = 633

# Row 63
$ws.Range("D63").Value = 44592
$ws.Range("I63").Value = "Tercera"
$ws.Range("J63").Value = 200
$ws.Range("K63").Value = 1500
$ws.Range("L63").Value = 1800
$ws.Range("M63").Value = 1650
$ws.Range("P63").Value = 550
